# CU11 - Compartir Reporte Inventario: requirement-text edits.
# Each change below splits (or merges) the w:t runs of a single
# paragraph. We use Range.InsertXML with a pkg:package payload that
# reproduces the exact original <w:p> attributes (w14:paraId,
# w14:textId, w:rsidR, w:rsidRDefault, w:rsidP) and <w:pPr> so that only
# the run content changes, matching the target OOXML precisely.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Insert-ParagraphXml([string]$findText, [string]$bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $findText"
    }
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml)
}

$rPr = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/></w:rPr>'
$pPr = '<w:pPr>' + $rPr + '</w:pPr>'

# 1) "El dueño podrá compartir el reporte de inventario al correo o
#    WhatsApp." -> split into two runs, dropping the "al correo o
#    WhatsApp." ending in favour of "o imprimir o descargar el mismo."
$body1 = '<w:p w14:paraId="4641BB76" w14:textId="601FBFA9" w:rsidR="00BF3851" w:rsidRDefault="00D02FA9">' + $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">El dueño podrá compartir el reporte de inventario </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>o imprimir o descargar el mismo.</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphXml "El dueño podrá compartir el reporte de inventario al correo o WhatsApp." $body1

# 2) "El dueño podrá enviar el reporte" -> split into three runs,
#    inserting " o descargar" before " el reporte".
$body2 = '<w:p w14:paraId="7796122E" w14:textId="69558B49" w:rsidR="007B30E6" w:rsidRDefault="00D02FA9" w:rsidP="007B30E6">' + $pPr + `
    '<w:r>' + $rPr + '<w:t>El dueño podrá enviar</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> o descargar</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> el reporte</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphXml "El dueño podrá enviar el reporte" $body2

# 3) "El sistema tendrá ... del reporte." -> drop trailing period, add
#    a second run ", o podrá descargar el mismo".
$body3 = '<w:p w14:paraId="534455F6" w14:textId="7C6D8B19" w:rsidR="007B30E6" w:rsidRDefault="00D02FA9" w:rsidP="007B30E6">' + $pPr + `
    '<w:r>' + $rPr + '<w:t>El sistema tendrá un formulario el cual le solicitará nombre, correo y numero para poder hacer el envío del reporte</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>, o podrá descargar el mismo</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphXml "El sistema tendrá un formulario el cual le solicitará nombre, correo y numero para poder hacer el envío del reporte." $body3

# 4) "<nº de veces> veces / <unidad de tiempo>" -> merge the three
#    runs (and the proofErr spell-check markers around "nº") into a
#    single run.
$body4 = '<w:p w14:paraId="35B6AD26" w14:textId="77777777" w:rsidR="007B30E6" w:rsidRDefault="007B30E6" w:rsidP="007B30E6">' + $pPr + `
    '<w:r>' + $rPr + '<w:t>&lt;nº de veces&gt; veces / &lt;unidad de tiempo&gt;</w:t></w:r>' + `
    '</w:p>'
Insert-ParagraphXml "<nº de veces> veces / <unidad de tiempo>" $body4

Write-Output "done"
